$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old (years-across-columns) table ------------------------
# The previous layout used B1:E1 as year headers and A2:A4 as metric
# labels, spanning A1:E4. Wipe contents AND the row-level formatting
# (row 1 had a custom row style) so nothing is left behind before the
# new, transposed table is written.
$ws.Range("A1:E4").ClearContents()
$ws.Rows("1:4").ClearFormats()

# --- New header row (metrics across columns, "year" label in A1) -------
$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = "% of U.S. adults who own smart speakers"
$ws.Range("C1").Value = "Number of smart speakers in US households"
$ws.Range("D1").Value = "Average number of smart speakers in U.S. households"

# Format the year column as text BEFORE typing the numeric-looking
# strings so Excel stores "2017"/"2018"/"2019" as text, matching the
# original file's text-formatted year values.
$ws.Range("A1:A4").NumberFormat = "@"
$ws.Range("A2").Value = "2017"
$ws.Range("A3").Value = "2018"
$ws.Range("A4").Value = "2019"

# --- Data rows, one per year --------------------------------------------
$ws.Range("B2").Value = 18
$ws.Range("C2").Value = 67000000
$ws.Range("D2").Value = 1.7

$ws.Range("B3").Value = 21
$ws.Range("C3").Value = 119000000
$ws.Range("D3").Value = 2.3

$ws.Range("B4").Value = 24
$ws.Range("C4").Value = 157000000
$ws.Range("D4").Value = 2.6

# "Number of smart speakers" column keeps its thousands-separator format.
$ws.Range("C2:C4").NumberFormat = "#,##0"

# Column E is no longer part of the table; make sure nothing is left there.
$ws.Range("E1:E4").Clear()

$ws.Range("D1").Select()
